$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 2 "Bitcoin"
Set-TextCell 2 3 "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
Set-TextCell 2 4 "69.940.59"
Set-TextCell 2 5 "  +0.15%  "

Set-TextCell 3 2 "Ethereum"
Set-TextCell 3 3 "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
Set-TextCell 3 4 "3.505.22"
Set-TextCell 3 5 "  -0.89%  "

Set-TextCell 4 2 "TetherUSD"
Set-TextCell 4 3 "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
Set-TextCell 4 4 "0.999"
Set-TextCell 4 5 "  -0.07%  "

Set-TextCell 5 2 "BNB"
Set-TextCell 5 3 "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextCell 5 4 "599.26"
Set-TextCell 5 5 "  -1.79%  "

Set-TextCell 6 2 "Solana"
Set-TextCell 6 3 "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell 6 4 "195.35"
Set-TextCell 6 5 "  +5.64%  "

Set-TextCell 7 2 "XRP"
Set-TextCell 7 3 "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextCell 7 4 "0.624"
Set-TextCell 7 5 "  +1.52%  "

Set-TextCell 8 2 "USDC"
Set-TextCell 8 3 "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextCell 8 4 "0.999"
Set-TextCell 8 5 "  -0.09%  "

Set-TextCell 9 2 "Dogecoin"
Set-TextCell 9 3 "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell 9 4 "0.209"
Set-TextCell 9 5 "  -1.36%  "

Set-TextCell 10 2 "Cardano"
Set-TextCell 10 3 "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell 10 4 "0.652"
Set-TextCell 10 5 "  +1.89%  "

Set-TextCell 11 2 "Avalanche"
Set-TextCell 11 3 "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell 11 4 "54.01"
Set-TextCell 11 5 "  +0.79%  "

Set-TextCell 12 2 "ShibaInu"
Set-TextCell 12 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell 12 4 "0.0000301"
Set-TextCell 12 5 "  -2.40%  "

Set-TextCell 13 2 "Polkadot"
Set-TextCell 13 3 "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell 13 4 "9.55"
Set-TextCell 13 5 "  +0.99%  "

Set-TextCell 14 2 "WrappedliquidstakedEther2.0"
Set-TextCell 14 3 "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell 14 4 "4.057.53"
Set-TextCell 14 5 "  -1.09%  "

Set-TextCell 15 2 "BitcoinCash"
Set-TextCell 15 3 "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell 15 4 "605.55"
Set-TextCell 15 5 "  +3.49%  "

Set-TextCell 16 2 "WrappedBTC"
Set-TextCell 16 3 "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell 16 4 "70.040.46"
Set-TextCell 16 5 "  +0.16%  "

Set-TextCell 17 2 "Chainlink"
Set-TextCell 17 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell 17 4 "18.98"
Set-TextCell 17 5 "  +0.60%  "

Set-TextCell 18 2 "Uniswap"
Set-TextCell 18 3 "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell 18 4 "12.61"
Set-TextCell 18 5 "  -0.24%  "

Set-TextCell 19 2 "WrappedEther"
Set-TextCell 19 3 "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell 19 4 "3.503.65"
Set-TextCell 19 5 "  -0.68%  "

Set-TextCell 20 2 "TRON"
Set-TextCell 20 3 "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell 20 4 "0.121"
Set-TextCell 20 5 "  +0.70%  "

Set-TextCell 21 2 "Polygon"
Set-TextCell 21 3 "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell 21 4 "0.992"
Set-TextCell 21 5 "  +0.12%  "

Set-TextCell 22 2 "InternetComputer(DFINITY)"
Set-TextCell 22 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 22 4 "18.15"
Set-TextCell 22 5 "  +3.78%  "

Set-TextCell 23 2 "Litecoin"
Set-TextCell 23 3 "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell 23 4 "104.51"
Set-TextCell 23 5 "  +9.01%  "

Set-TextCell 24 2 "Toncoin"
Set-TextCell 24 3 "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell 24 4 "5.12"
Set-TextCell 24 5 "  +5.74%  "

Set-TextCell 25 2 "PancakeSwap"
Set-TextCell 25 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell 25 4 "4.58"
Set-TextCell 25 5 "  -2.56%  "

Set-TextCell 26 2 "ImmutableX"
Set-TextCell 26 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell 26 4 "3.06"
Set-TextCell 26 5 "  +3.04%  "

Set-TextCell 27 2 "RenderToken"
Set-TextCell 27 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 27 4 "10.94"
Set-TextCell 27 5 "  -0.39%  "

Set-TextCell 28 2 "Filecoin"
Set-TextCell 28 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell 28 4 "9.72"
Set-TextCell 28 5 "  +1.87%  "

Set-TextCell 29 2 "EthereumClassic"
Set-TextCell 29 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell 29 4 "33.55"
Set-TextCell 29 5 "  +4.83%  "

Set-TextCell 30 2 "dogwifhat"
Set-TextCell 30 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell 30 4 "4.57"
Set-TextCell 30 5 "  +26.56%  "

Set-TextCell 31 2 "NEARProtocol"
Set-TextCell 31 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 31 4 "7.12"
Set-TextCell 31 5 "  +1.75%  "

Set-TextCell 32 2 "Cosmos"
Set-TextCell 32 3 "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell 32 4 "12.64"
Set-TextCell 32 5 "  +4.11%  "

Set-TextCell 33 2 "Hedera"
Set-TextCell 33 3 "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell 33 4 "0.115"
Set-TextCell 33 5 "  +1.39%  "

Set-TextCell 34 2 "OKB"
Set-TextCell 34 3 "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell 34 4 "63.45"
Set-TextCell 34 5 "  +0.27%  "

Set-TextCell 35 2 "PEPE"
Set-TextCell 35 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell 35 4 "0.0₃0824"
Set-TextCell 35 5 "  +6.14%  "

Set-TextCell 36 2 "Maker"
Set-TextCell 36 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell 36 4 "3.726.56"
Set-TextCell 36 5 "  +5.49%  "

Set-TextCell 37 2 "Dai"
Set-TextCell 37 3 "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell 37 4 "1.00"
Set-TextCell 37 5 "  -0.09%  "

Set-TextCell 38 2 "Fetch.AI"
Set-TextCell 38 3 "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell 38 4 "3.05"
Set-TextCell 38 5 "  -6.64%  "

Set-TextCell 39 2 "TheGraph"
Set-TextCell 39 3 "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell 39 4 "0.392"
Set-TextCell 39 5 "  -2.27%  "

Set-TextCell 40 2 "InjectiveProtocol"
Set-TextCell 40 3 "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell 40 4 "36.74"
Set-TextCell 40 5 "  -0.96%  "

Set-TextCell 41 2 "Stacks"
Set-TextCell 41 3 "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextCell 41 4 "3.57"
Set-TextCell 41 5 "  +2.29%  "

Set-TextCell 42 2 "Bittensor"
Set-TextCell 42 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell 42 4 "499.07"
Set-TextCell 42 5 "  -5.79%  "

Set-TextCell 43 2 "Kaspa"
Set-TextCell 43 3 "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell 43 4 "0.136"
Set-TextCell 43 5 "  +0.08%  "

Set-TextCell 44 2 "VeChain"
Set-TextCell 44 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 44 4 "0.0457"
Set-TextCell 44 5 "  +0.59%  "

Set-TextCell 45 2 "Stellar"
Set-TextCell 45 3 "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell 45 4 "0.140"
Set-TextCell 45 5 "  -1.20%  "

Set-TextCell 46 2 "ThetaToken"
Set-TextCell 46 3 "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell 46 4 "2.83"
Set-TextCell 46 5 "  -3.10%  "

Set-TextCell 47 2 "ApeXProtocol"
Set-TextCell 47 3 "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextCell 47 4 "3.32"
Set-TextCell 47 5 "  -1.33%  "

Set-TextCell 48 2 "FirstDigitalUSD"
Set-TextCell 48 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell 48 4 "1.01"
Set-TextCell 48 5 "  +0.39%  "

Set-TextCell 49 2 "THORChain"
Set-TextCell 49 3 "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextCell 49 4 "8.71"
Set-TextCell 49 5 "  -4.70%  "

Set-TextCell 50 2 "FLOKI"
Set-TextCell 50 3 "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextCell 50 4 "0.000245"
Set-TextCell 50 5 "  +1.97%  "

Set-TextCell 51 2 "Monero"
Set-TextCell 51 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 51 4 "131.62"
Set-TextCell 51 5 "  -2.86%  "

Write-Output "Updated cryptos list"
